# Update cryptos list - GitHub Actions scheduled data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write a value as plain text, preventing Excel from auto-coercing
# numeric-looking strings (e.g. "202.14") into Double values, and
# without leaving a residual NumberFormat/style on the cell.
function Set-Text($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

function Set-Row($Row, $B, $C, $D, $E) {
    if ($B) { Set-Text "B$Row" $B }
    if ($C) { Set-Text "C$Row" $C }
    if ($D) { Set-Text "D$Row" $D }
    if ($E) { Set-Text "E$Row" $E }
}

Set-Row 2  $null $null "75.981.55" "  +1.48%  "
Set-Row 3  $null $null "2.930.83"  "  +4.49%  "
Set-Row 4  $null $null $null       "  +0.08%  "
Set-Row 5  $null $null "202.14"    "  +7.71%  "
Set-Row 6  $null $null "597.73"    "  +1.14%  "
Set-Row 7  $null $null "1.00"      "  +0.04%  "
Set-Row 8  $null $null "0.554"     "  +1.97%  "
Set-Row 9  $null $null "0.197"     "  +4.90%  "
Set-Row 10 $null $null "2.938.14"  "  +4.87%  "
Set-Row 11 $null $null "0.448"     "  +19.94%  "
Set-Row 12 $null $null $null       "  +1.03%  "
Set-Row 13 $null $null $null       "  +2.28%  "
Set-Row 14 $null $null "3.473.26"  "  +4.67%  "
Set-Row 15 $null $null "28.39"     "  +5.87%  "
Set-Row 16 $null $null "75.944.81" "  +1.62%  "
Set-Row 17 $null $null "0.0000190" "  +2.58%  "
Set-Row 18 $null $null "2.943.73"  "  +4.78%  "
Set-Row 19 $null $null "13.27"     "  +8.64%  "
Set-Row 20 $null $null "8.79"      "  -0.70%  "
Set-Row 21 $null $null "374.86"    "  -0.17%  "
Set-Row 22 $null $null "2.32"      "  +2.58%  "
Set-Row 23 $null $null "4.36"      "  +7.00%  "
Set-Row 24 $null $null $null       "  +1.71%  "
Set-Row 25 $null $null $null       "  -0.01%  "

# Rows 26/27 swap (NEARProtocol <-> WrappedeETH)
Set-Row 26 "WrappedeETH"  "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"  "3.091.47" "  +4.78%  "
Set-Row 27 "NEARProtocol" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"  "4.35"     "  +5.40%  "

Set-Row 28 $null $null "9.73"      "  +1.16%  "
Set-Row 29 $null $null "0.0000109" "  +6.77%  "
Set-Row 30 $null $null "1.00"      "  +0.02%  "
Set-Row 31 $null $null "1.39"      "  +0.23%  "
Set-Row 32 $null $null "7.96"      "  +5.17%  "
Set-Row 33 $null $null "501.74"    "  -1.34%  "
Set-Row 34 $null $null "1.85"      "  +3.79%  "
Set-Row 35 $null $null "1.00"      "  +0.17%  "

# Rows 36/37 swap (Monero <-> EthereumClassic)
Set-Row 36 "EthereumClassic" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc" "20.34"  "  +2.88%  "
Set-Row 37 "Monero"          "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"          "163.94" "  -0.26%  "

Set-Row 38 $null $null $null    "  +27.07%  "
Set-Row 39 $null $null "19.64"  "  +1.42%  "
Set-Row 40 $null $null "0.374"  "  +10.37%  "
Set-Row 41 $null $null $null    "  -3.40%  "
Set-Row 42 $null $null $null    "  +0.03%  "
Set-Row 43 $null $null "179.61" "  -0.69%  "
Set-Row 44 $null $null "5.01"   "  +0.89%  "
Set-Row 45 $null $null "1.67"   "  +1.12%  "

# Rows 46/47 swap (ImmutableX <-> OKB)
Set-Row 46 "OKB"        "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"     "40.17" "  +0.44%  "
Set-Row 47 "ImmutableX"  "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx" "1.20"  "  +0.63%  "

# Rows 48/49 swap (dogwifhat <-> Filecoin)
Set-Row 48 "Filecoin"   "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"   "3.95" "  +6.68%  "
Set-Row 49 "dogwifhat"  "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"  "2.34" "  +2.04%  "

Set-Row 50 $null $null "0.581" "  +2.28%  "
Set-Row 51 $null $null "22.77" "  +9.93%  "

$wb.Save()
